## Fruta / hortaliza, semanal
## Insert 3 new weekly-report rows (Palta, Terminal Hortofrutícola Agro
## Chillán) right before the current row 1192, pushing the existing
## rows 1192:1257 down to 1195:1260, then populate the 3 new rows with
## the latest week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above row 1192; Excel shifts row 1192:1257
# down to 1195:1260 and the dimension grows to A1:T1260.
$ws.Rows("1192:1194").Insert()

# Shared (unchanged) categorical columns for all three new rows.
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$tipo        = "Fruta"
$productoId  = 100106
$producto    = "Oleaginosos"
$categoriaId = 100106002
$categoria   = "Palta"
$variedad    = "Hass"
$origen      = "Provincia de Quillota"
$fecha       = 45267

# Row 1192 - Especial
$r = 1192
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 150
$ws.Cells.Item($r, 14).Value = 3500
$ws.Cells.Item($r, 15).Value = 3500
$ws.Cells.Item($r, 16).Value = 3500
$ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 3500
$ws.Cells.Item($r, 20).Value = 1

# Row 1193 - Primera
$r = 1193
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 150
$ws.Cells.Item($r, 14).Value = 3000
$ws.Cells.Item($r, 15).Value = 3000
$ws.Cells.Item($r, 16).Value = 3000
$ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 3000
$ws.Cells.Item($r, 20).Value = 1

# Row 1194 - Segunda
$r = 1194
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 150
$ws.Cells.Item($r, 14).Value = 2800
$ws.Cells.Item($r, 15).Value = 2800
$ws.Cells.Item($r, 16).Value = 2800
$ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2800
$ws.Cells.Item($r, 20).Value = 1
